$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-1.238285183906555, 0.4941155612468719, -1.145785808563232, -0.1621847152709961, -0.0795651972293853, -0.5577199459075928),
    @(0.3425590991973889, 0.401696681976318, -1.307540893554688, -0.1879937797784805, -0.1305724531412124, -0.0817032232880592),
    @(1.148096084594726, 0.2248815298080447, -1.84254863858223, -0.0765108689665794, -0.2229658216238021, -0.0765108689665794),
    @(0.3787193298339841, 0.750096321105957, -2.375997304916381, -0.0930042341351509, 0.1044579595327377, -0.295353353023529),
    @(0.1436797380447386, 0.7403357923030853, -1.325827866792679, 0.0355829000473022, 0.3640756905078888, -0.0638354197144508),
    @(-1.716686248779299, 1.053612291812897, 0.4070562124252333, 0.5752823352813721, 0.5053382515907288, -0.007177666760981),
    @(-3.969705402851109, 0.09983259439468162, 0.9274015724658967, 0.5288565754890442, 0.7619016766548157, 0.2580905556678772),
    @(-5.239097833633425, -0.5959589481353762, -0.1143757104873694, -0.7357872128486633, -0.0797179117798805, 0.8136724829673767),
    @(2.036354780197156, -0.2977316975593562, -1.156379550695419, -1.146135926246643, -0.6188064813613892, -0.1640173196792602),
    @(5.741946458816498, -4.677600264549274, -8.871290028095274, 0.8213083148002625, -3.046079635620117, 1.09803032875061),
    @(-18.10554087162019, -17.49406802654266, -29.53683829307556, 1.518305540084839, -0.5958990454673767, 0.1006400510668754),
    @(1.688319206237781, -2.673514366149895, -12.34325218200682, -0.9292787313461304, 0.0415388382971286, -2.335643291473389),
    @(-3.507262408733383, 0.2684899270534595, -2.699394106864903, 0.5500841736793518, 1.703092336654663, -0.4928155243396759),
    @(-0.715216994285579, 0.4368197321891781, 0.5279676914215112, 0.2722931802272796, -0.3020728528499603, 0.2756529450416565),
    @(0.77803122997284, 0.2878375947475432, -0.8067402243614243, -0.6982190012931824, 0.4492913782596588, -0.5900958180427551),
    @(-0.02435183525085627, -0.3931519985198989, -2.55875074863434, -0.4308127164840698, 1.22447943687439, -0.319024384021759),
    @(-0.2728092074394176, -0.8767854124307625, -2.626779749989506, -0.0684169083833694, 0.4977024495601654, -0.0363464802503585),
    @(1.066039085388185, -0.5432969331741325, -1.81298840045929, -0.1298088580369949, 0.0172569435089826, 0.1985312104225158),
    @(1.15899240970611, -0.2723855525255186, -1.873978555202483, 0.0598647929728031, 0.5253441333770752, -0.0415388382971286),
    @(-0.3004360198974618, 0.163273096084592, -1.427715420722963, -0.1009454801678657, 0.2724458873271942, -0.4051563739776611),
    @(-0.5053797960281373, -0.4524855315685276, -1.996987149119378, -0.0010690141934901, -0.1505782902240753, 0.0145080499351024),
    @(0.1984210014343271, -0.3136261999607083, -1.604727536439895, -0.1108720451593399, 0.0652098655700683, -0.131183311343193),
    @(0.1348390579223627, -0.3085210472345353, -1.704802349209786, 0.0563523173332214, 0.0316122770309448, 0.1798998117446899),
    @(-0.11956262588501, -0.3192775845527647, -1.924065947532654, -0.0160352122038602, 0.0537561401724815, -0.0355829000473022),
    @(-0.04362952709197906, -0.3616583049297337, -1.960778713226318, -0.0247400421649217, 0.0786489024758338, -0.0580321997404098),
    @(0.191988468170166, -0.4578718543052673, -1.877070605754853, -0.026419922709465, 0.0665843114256858, -0.0759000033140182),
    @(0.1683353185653686, -0.3762182295322418, -1.764573842287064, 0.0235183127224445, 0.012980886735022, -0.0116064399480819),
    @(0.1497325897216795, -0.3715704679489131, -1.762347698211669, -0.0093156946823, 0.0103847095742821, -0.0035124751739203),
    @(0.0995370149612424, -0.2668604403734202, -1.60942207276821, -0.0279470849782228, 0.0250454749912023, -0.00335975876078),
    @(0.02663779258727977, -0.2185956239700315, -1.59472194314003, -0.0271835029125213, 0.0594066455960273, -0.0343611687421798)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = ($i * 100)
    $ws.Cells.Item($r, 2).Value = "falling"
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, 3 + $j).Value = $row[$j]
    }
}
